$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.296.06"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "3.679.19"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'675.18"
$ws.Range("D6").Value = "'158.28"
$ws.Range("E6").Value = "  -2.46%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -1.27%  "
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("E10").Value = "  -6.19%  "
$ws.Range("E11").Value = "  -2.27%  "
$ws.Range("E12").Value = "  -3.72%  "
$ws.Range("D13").Value = "4.299.52"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("E14").Value = "  -3.89%  "
$ws.Range("D15").Value = "3.675.95"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").Value = "69.257.42"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("E17").Value = "  +1.62%  "
$ws.Range("E18").Value = "  -2.05%  "
$ws.Range("E19").Value = "  -3.07%  "
$ws.Range("D20").Value = "'468.80"
$ws.Range("E20").Value = "  -2.88%  "
$ws.Range("D21").Value = "'9.94"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("E22").Value = "  -2.39%  "
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("D24").Value = "3.822.06"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  -5.98%  "
$ws.Range("E28").Value = "  -4.63%  "
$ws.Range("E29").Value = "  -1.69%  "
$ws.Range("E30").Value = "  -4.33%  "
$ws.Range("D31").Value = "'6.62"
$ws.Range("E31").Value = "  -3.22%  "
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'26.94"
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("E34").Value = "  -4.86%  "
$ws.Range("D35").Value = "3.670.78"
$ws.Range("E36").Value = "  -5.04%  "
$ws.Range("E37").Value = "  -3.38%  "
$ws.Range("E38").Value = "  -1.37%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").Value = "'2.22"
$ws.Range("E41").Value = "  -4.43%  "
$ws.Range("D42").Value = "'0.0903"
$ws.Range("E42").Value = "  -3.52%  "
$ws.Range("D43").Value = "'170.88"
$ws.Range("E43").Value = "  +4.74%  "
$ws.Range("D44").Value = "'0.942"
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("E45").Value = "  -1.54%  "
$ws.Range("B46").Value = "FLOKI"
$ws.Range("C46").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D46").Value = "'0.000276"
$ws.Range("E46").Value = "  -3.95%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "'2.68"
$ws.Range("E47").Value = "  -5.47%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'27.76"
$ws.Range("E48").Value = "  -7.25%  "
$ws.Range("E49").Value = "  -3.64%  "
$ws.Range("D50").Value = "'1.09"
$ws.Range("E50").Value = "  -3.19%  "
$ws.Range("E51").Value = "  -3.16%  "
